$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Gnai2"
$ws.Range("C2").Value = "F2r"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 162.399297
$ws.Range("H2").Value = 487.197891
$ws.Range("I2").Value = 0.3910371682630009
$ws.Range("J2").Value = 0.3910371682630009
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 6.153886666666666
$ws.Range("N2").Value = 18.46166
$ws.Range("O2").Value = 0.08077417226496708
$ws.Range("P2").Value = 0.0807741722649671
$ws.Range("Q2").Value = 999.3868684843399
$ws.Range("R2").Value = 8994.481816359059
$ws.Range("S2").Value = 0.03158570359128056
$ws.Range("T2").Value = 0.03158570359128056

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Gnai2"
$ws.Range("C3").Value = "F2r"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 162.399297
$ws.Range("H3").Value = 487.197891
$ws.Range("I3").Value = 0.3910371682630009
$ws.Range("J3").Value = 0.3910371682630009
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 48.59929933333333
$ws.Range("N3").Value = 145.797898
$ws.Range("O3").Value = 0.6379006291374719
$ws.Range("P3").Value = 0.637900629137472
$ws.Range("Q3").Value = 7892.492046425902
$ws.Range("R3").Value = 71032.42841783311
$ws.Range("S3").Value = 0.2494428556511038
$ws.Range("T3").Value = 0.2494428556511038

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Gnai2"
$ws.Range("C4").Value = "F2r"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 162.399297
$ws.Range("H4").Value = 487.197891
$ws.Range("I4").Value = 0.3910371682630009
$ws.Range("J4").Value = 0.3910371682630009
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 21.43313066666667
$ws.Range("N4").Value = 64.299392
$ws.Range("O4").Value = 0.2813251985975609
$ws.Range("P4").Value = 0.2813251985975609
$ws.Range("Q4").Value = 3480.725352775808
$ws.Range("R4").Value = 31326.52817498227
$ws.Range("S4").Value = 0.1100086090206166
$ws.Range("T4").Value = 0.1100086090206166

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Gnai2"
$ws.Range("C5").Value = "F2r"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 65.41736466666667
$ws.Range("H5").Value = 196.252094
$ws.Range("I5").Value = 0.1575168212364948
$ws.Range("J5").Value = 0.1575168212364948
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 6.153886666666666
$ws.Range("N5").Value = 18.46166
$ws.Range("O5").Value = 0.08077417226496708
$ws.Range("P5").Value = 0.0807741722649671
$ws.Range("Q5").Value = 402.5710481906711
$ws.Range("R5").Value = 3623.13943371604
$ws.Range("S5").Value = 0.01272329085318666
$ws.Range("T5").Value = 0.01272329085318666

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Gnai2"
$ws.Range("C6").Value = "F2r"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 65.41736466666667
$ws.Range("H6").Value = 196.252094
$ws.Range("I6").Value = 0.1575168212364948
$ws.Range("J6").Value = 0.1575168212364948
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 48.59929933333333
$ws.Range("N6").Value = 145.797898
$ws.Range("O6").Value = 0.6379006291374719
$ws.Range("P6").Value = 0.637900629137472
$ws.Range("Q6").Value = 3179.238087033157
$ws.Range("R6").Value = 28613.14278329841
$ws.Range("S6").Value = 0.1004800793664948
$ws.Range("T6").Value = 0.1004800793664948

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Gnai2"
$ws.Range("C7").Value = "F2r"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 65.41736466666667
$ws.Range("H7").Value = 196.252094
$ws.Range("I7").Value = 0.1575168212364948
$ws.Range("J7").Value = 0.1575168212364948
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 21.43313066666667
$ws.Range("N7").Value = 64.299392
$ws.Range("O7").Value = 0.2813251985975609
$ws.Range("P7").Value = 0.2813251985975609
$ws.Range("Q7").Value = 1402.09892476965
$ws.Range("R7").Value = 12618.89032292685
$ws.Range("S7").Value = 0.04431345101681341
$ws.Range("T7").Value = 0.04431345101681341

$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Gnai2"
$ws.Range("C8").Value = "F2r"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 126.3069433333333
$ws.Range("H8").Value = 378.92083
$ws.Range("I8").Value = 0.3041313008456065
$ws.Range("J8").Value = 0.3041313008456065
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 6.153886666666666
$ws.Range("N8").Value = 18.46166
$ws.Range("O8").Value = 0.08077417226496708
$ws.Range("P8").Value = 0.0807741722649671
$ws.Range("Q8").Value = 777.2786144864222
$ws.Range("R8").Value = 6995.507530377799
$ws.Range("S8").Value = 0.02456595408567154
$ws.Range("T8").Value = 0.02456595408567155

$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Gnai2"
$ws.Range("C9").Value = "F2r"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 126.3069433333333
$ws.Range("H9").Value = 378.92083
$ws.Range("I9").Value = 0.3041313008456065
$ws.Range("J9").Value = 0.3041313008456065
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 48.59929933333333
$ws.Range("N9").Value = 145.797898
$ws.Range("O9").Value = 0.6379006291374719
$ws.Range("P9").Value = 0.637900629137472
$ws.Range("Q9").Value = 6138.428946935038
$ws.Range("R9").Value = 55245.86052241534
$ws.Range("S9").Value = 0.1940055481498101
$ws.Range("T9").Value = 0.1940055481498102

$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Gnai2"
$ws.Range("C10").Value = "F2r"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 126.3069433333333
$ws.Range("H10").Value = 378.92083
$ws.Range("I10").Value = 0.3041313008456065
$ws.Range("J10").Value = 0.3041313008456065
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 21.43313066666667
$ws.Range("N10").Value = 64.299392
$ws.Range("O10").Value = 0.2813251985975609
$ws.Range("P10").Value = 0.2813251985975609
$ws.Range("Q10").Value = 2707.153220570596
$ws.Range("R10").Value = 24364.37898513536
$ws.Range("S10").Value = 0.08555979861012479
$ws.Range("T10").Value = 0.0855597986101248

$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Gnai2"
$ws.Range("C11").Value = "F2r"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 61.180387
$ws.Range("H11").Value = 183.541161
$ws.Range("I11").Value = 0.1473147096548978
$ws.Range("J11").Value = 0.1473147096548978
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 6.153886666666666
$ws.Range("N11").Value = 18.46166
$ws.Range("O11").Value = 0.08077417226496708
$ws.Range("P11").Value = 0.0807741722649671
$ws.Range("Q11").Value = 376.4971678208066
$ws.Range("R11").Value = 3388.47451038726
$ws.Range("S11").Value = 0.01189922373482832
$ws.Range("T11").Value = 0.01189922373482833

$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Gnai2"
$ws.Range("C12").Value = "F2r"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 61.180387
$ws.Range("H12").Value = 183.541161
$ws.Range("I12").Value = 0.1473147096548978
$ws.Range("J12").Value = 0.1473147096548978
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 48.59929933333333
$ws.Range("N12").Value = 145.797898
$ws.Range("O12").Value = 0.6379006291374719
$ws.Range("P12").Value = 0.637900629137472
$ws.Range("Q12").Value = 2973.323941142175
$ws.Range("R12").Value = 26759.91547027958
$ws.Range("S12").Value = 0.09397214597006331
$ws.Range("T12").Value = 0.09397214597006334

$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Gnai2"
$ws.Range("C13").Value = "F2r"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 61.180387
$ws.Range("H13").Value = 183.541161
$ws.Range("I13").Value = 0.1473147096548978
$ws.Range("J13").Value = 0.1473147096548978
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 21.43313066666667
$ws.Range("N13").Value = 64.299392
$ws.Range("O13").Value = 0.2813251985975609
$ws.Range("P13").Value = 0.2813251985975609
$ws.Range("Q13").Value = 1311.287228808235
$ws.Range("R13").Value = 11801.58505927411
$ws.Range("S13").Value = 0.04144333995000615
$ws.Range("T13").Value = 0.04144333995000616
